# Regenerate save_data column G ("K") with freshly computed values.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")
#
# Column G holds the "K" metric for each logged game/row (rows 2-74).
# The values below are the newly-calculated replacements for that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 3
    9  = 3
    10 = 1
    11 = 3
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 5
    19 = 3
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 3
    29 = 3
    30 = 0
    31 = 4
    32 = 3
    33 = 1
    34 = 2
    35 = 3
    36 = 1
    37 = 2
    38 = 0
    39 = 4
    40 = 1
    41 = 3
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 2
    47 = 2
    48 = 1
    49 = 0
    50 = 0
    51 = 1
    52 = 0
    53 = 1
    54 = 0
    55 = 1
    56 = 2
    57 = 1
    58 = 2
    59 = 1
    60 = 0
    61 = 1
    62 = 1
    63 = 1
    64 = 4
    65 = 0
    66 = 2
    67 = 1
    68 = 2
    69 = 2
    70 = 3
    71 = 1
    72 = 2
    73 = 2
    74 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
